# Insert a new student record ("Bintang Diyantoro") as the first data row
# of the roster, pushing every existing row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right below the header (row 1), shifting all
# existing student rows (old rows 2-13) down to rows 3-14.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new student's data.
$ws.Range("A2").Value = "Bintang Diyantoro"
$ws.Range("B2").Value = 12345678
$ws.Range("C2").Value = "bintangdiyantoro@gmail.com"
$ws.Range("D2").Value = "Technique Informatique"
